$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.650.50"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.596.82"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.19%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'211.31"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.67%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.18%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.47%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -0.97%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0837"
$ws.Range("E11").Value = "  -0.37%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.821.85"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13 - was Polkadot, now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.591.29"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14 - was WrappedEther, now Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.02"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.12%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.11"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.647.25"
$ws.Range("E17").Value = "  -0.16%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.24%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'209.48"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.15%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +3.97%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.46%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +0.75%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'145.01"
$ws.Range("E25").Value = "  -1.09%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -1.11%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.03%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.36%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  +0.56%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.43%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.281.57"
$ws.Range("E34").Value = "  -0.96%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -7.69%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.64%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +1.16%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.87%  "

# Row 39 - WEMIXToken
$ws.Range("D39").Value = "'1.08"
$ws.Range("E39").Value = "  +22.31%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.08%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "'5.53"
$ws.Range("E41").Value = "  +2.79%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -0.01%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'0.785"
$ws.Range("E43").Value = "  -0.67%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'63.92"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.733.89"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46 - Quant
$ws.Range("D46").Value = "'90.86"
$ws.Range("E46").Value = "  +0.77%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -3.55%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +2.34%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.58%  "

# Row 50 - USDD
$ws.Range("E50").Value = "  -0.16%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "'7.38"
$ws.Range("E51").Value = "  -1.64%  "
